# IP_setting/saved_adresses_2.xlsx - "nová verze ip setting, oprava chyb"
#
# Updates the ip_adress_list sheet contents (rows 1-4 edited, row 5 added),
# and moves the active-cell selection on both the ip_adress_list sheet and
# the Settings sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ip_adress_list")
$ws3 = $wb.Worksheets.Item("Settings")

# ---------------------------------------------------------------------
# Row 1
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "sgf"
$ws1.Range("B1").Value = "192.168.100.241"
$ws1.Range("C1").Value = "255.255.255.0"
$ws1.Range("D1").Value = "sfg"

# ---------------------------------------------------------------------
# Row 2  (A2 becomes the numeric-looking text "514" - force text type by
# writing a non-numeric placeholder with a proper text number format and
# then restoring the Normal style so no stray format is left behind)
# ---------------------------------------------------------------------
$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "514"
$ws1.Range("A2").Style = "Normal"
$ws1.Range("B2").Value = "192.168.100.241"
$ws1.Range("C2").Value = "255.255.255.0"
$ws1.Range("D2").Value = "afs`nasdf`nasdf"
# The multi-line text auto-expands the row height; AutoFit restores the
# sheet's standard (non-custom) row height so no stray ht/customHeight is
# left on row 2 (matches the source diff, which leaves row 2 untouched).
$ws1.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = "sfdgsfg"
$ws1.Range("B3").Value = "192.168.100.241"
$ws1.Range("C3").Value = "255.255.255.0"
$ws1.Range("D3").Value = "sfs"

# ---------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------
$ws1.Range("A4").Value = "dsf"
$ws1.Range("B4").Value = "192.168.100.241"
$ws1.Range("C4").Value = "255.255.255.0"
$ws1.Range("D4").Value = "fffffffffffffffffffffffffffffff"

# ---------------------------------------------------------------------
# Row 5 (new)
# ---------------------------------------------------------------------
$ws1.Range("A5").Value = "sf"
$ws1.Range("B5").Value = "192.168.100.241"
$ws1.Range("C5").Value = "255.255.255.0"
$ws1.Range("D5").Value = "sdfsfdddddddddddddd"

# ---------------------------------------------------------------------
# Selections - ip_adress_list moves to F15, Settings moves to B14.
# Settings stays the active/visible tab (as before), so re-activate it
# after touching the ip_adress_list selection.
# ---------------------------------------------------------------------
$ws1.Range("F15").Select() | Out-Null
$ws3.Activate() | Out-Null
$ws3.Range("B14").Select() | Out-Null
